$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.301437
$ws.Range("H2").Value = 0.9043110000000001
$ws.Range("I2").Value = 0.003123224212368042
$ws.Range("J2").Value = 0.003123224212368043
$ws.Range("M2").Value = 9.084137666666667
$ws.Range("N2").Value = 27.252413
$ws.Range("O2").Value = 0.2765376761551382
$ws.Range("P2").Value = 0.2765376761551382
$ws.Range("Q2").Value = 2.738295205827
$ws.Range("R2").Value = 24.644656852443
$ws.Range("S2").Value = 0.0008636891657997203
$ws.Range("T2").Value = 0.0008636891657997204

$ws.Range("G3").Value = 0.301437
$ws.Range("H3").Value = 0.9043110000000001
$ws.Range("I3").Value = 0.003123224212368042
$ws.Range("J3").Value = 0.003123224212368043
$ws.Range("O3").Value = 0.3707916163717078
$ws.Range("P3").Value = 0.3707916163717078
$ws.Range("Q3").Value = 3.671604244269
$ws.Range("R3").Value = 33.04443819842101
$ws.Range("S3").Value = 0.001158065353995201
$ws.Range("T3").Value = 0.001158065353995201

$ws.Range("G4").Value = 0.301437
$ws.Range("H4").Value = 0.9043110000000001
$ws.Range("I4").Value = 0.003123224212368042
$ws.Range("J4").Value = 0.003123224212368043
$ws.Range("M4").Value = 11.58507333333333
$ws.Range("N4").Value = 34.75522
$ws.Range("O4").Value = 0.3526707074731541
$ws.Range("P4").Value = 0.3526707074731541
$ws.Range("Q4").Value = 3.49216975038
$ws.Range("R4").Value = 31.42952775342
$ws.Range("S4").Value = 0.001101469692573122
$ws.Range("T4").Value = 0.001101469692573122

$ws.Range("I5").Value = 0.7782793322359159
$ws.Range("J5").Value = 0.7782793322359159
$ws.Range("M5").Value = 9.084137666666667
$ws.Range("N5").Value = 27.252413
$ws.Range("O5").Value = 0.2765376761551382
$ws.Range("P5").Value = 0.2765376761551382
$ws.Range("Q5").Value = 682.3584921685766
$ws.Range("R5").Value = 6141.22642951719
$ws.Range("S5").Value = 0.2152235579360929
$ws.Range("T5").Value = 0.2152235579360929

$ws.Range("I6").Value = 0.7782793322359159
$ws.Range("J6").Value = 0.7782793322359159
$ws.Range("O6").Value = 0.3707916163717078
$ws.Range("P6").Value = 0.3707916163717078
$ws.Range("S6").Value = 0.2885794515884487
$ws.Range("T6").Value = 0.2885794515884487

$ws.Range("I7").Value = 0.7782793322359159
$ws.Range("J7").Value = 0.7782793322359159
$ws.Range("M7").Value = 11.58507333333333
$ws.Range("N7").Value = 34.75522
$ws.Range("O7").Value = 0.3526707074731541
$ws.Range("P7").Value = 0.3526707074731541
$ws.Range("Q7").Value = 870.2172359631845
$ws.Range("R7").Value = 7831.955123668661
$ws.Range("S7").Value = 0.2744763227113744
$ws.Range("T7").Value = 0.2744763227113744

$ws.Range("G8").Value = 21.09786333333333
$ws.Range("H8").Value = 63.29359
$ws.Range("I8").Value = 0.2185974435517159
$ws.Range("J8").Value = 0.218597443551716
$ws.Range("M8").Value = 9.084137666666667
$ws.Range("N8").Value = 27.252413
$ws.Range("O8").Value = 0.2765376761551382
$ws.Range("P8").Value = 0.2765376761551382
$ws.Range("Q8").Value = 191.6558949925189
$ws.Range("R8").Value = 1724.90305493267
$ws.Range("S8").Value = 0.06045042905324552
$ws.Range("T8").Value = 0.06045042905324553

$ws.Range("G9").Value = 21.09786333333333
$ws.Range("H9").Value = 63.29359
$ws.Range("I9").Value = 0.2185974435517159
$ws.Range("J9").Value = 0.218597443551716
$ws.Range("O9").Value = 0.3707916163717078
$ws.Range("P9").Value = 0.3707916163717078
$ws.Range("Q9").Value = 256.9790853799433
$ws.Range("R9").Value = 2312.81176841949
$ws.Range("S9").Value = 0.08105409942926391
$ws.Range("T9").Value = 0.08105409942926392

$ws.Range("G10").Value = 21.09786333333333
$ws.Range("H10").Value = 63.29359
$ws.Range("I10").Value = 0.2185974435517159
$ws.Range("J10").Value = 0.218597443551716
$ws.Range("M10").Value = 11.58507333333333
$ws.Range("N10").Value = 34.75522
$ws.Range("O10").Value = 0.3526707074731541
$ws.Range("P10").Value = 0.3526707074731541
$ws.Range("Q10").Value = 244.4202938933111
$ws.Range("R10").Value = 2199.7826450398
$ws.Range("S10").Value = 0.07709291506920653
$ws.Range("T10").Value = 0.07709291506920653
